$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Utilisateurs")

$ws.Range("D2").Value = '["jeudi","samedi","dimanche"]'
$ws.Range("F2").Value = '[{"transactionId":"1745665348518","payer":"Mila","amount":185,"description":"test","repayments":[{"userId":"Mila","amount":185,"paid":false}],"paid":false}]'
$ws.Range("F3").Value = '[{"transactionId":"1745665348518","payer":"Mila","amount":185,"description":"test","repayments":[{"userId":"Ju","amount":185,"paid":false}],"paid":false}]'
$ws.Range("F4").Value = '[{"transactionId":"1745665348518","payer":"Mila","amount":185,"description":"test","repayments":[{"userId":"Louise","amount":185,"paid":false}],"paid":false}]'
